$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add the new row for "E-Packet DDP" / "EPSDP" at row 18
$ws.Range("A18").Value = 98
$ws.Range("B18").Value = "E-Packet DDP"
$ws.Range("C18").Value = "EPSDP"

# Select the newly added cell to match the saved selection state
$ws.Range("B18").Select()
